$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '87.835.50'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.172.93'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.69%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '207.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '611.20'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.82%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.389'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.64%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.675'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.21%  '
$ws.Range('E9').Value = '  -0.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.167.57'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.538'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -8.95%  '
$ws.Range('E12').Value = '  -0.77%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000245'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -8.14%  '
$ws.Range('B14').Value = 'Toncoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.28'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.44%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.757.06'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.87%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '87.610.67'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '32.25'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -6.89%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.165.17'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.20'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.52%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.49'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -5.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '413.52'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.56%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.48'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -8.21%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.07'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -6.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.25'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.23'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.19%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.341.87'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.65%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000133'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '73.52'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.83%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.08%  '
$ws.Range('E30').Value = '  -11.11%  '
$ws.Range('E31').Value = '  -0.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '547.02'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.58%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '8.24'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -8.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.33'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -8.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.93'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.93%  '
$ws.Range('E36').Value = '  -6.27%  '
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '21.92'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.83%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.130'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.97%  '
$ws.Range('E39').Value = '  +0.30%  '
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('E41').Value = '  +1.57%  '
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('E43').Value = '  -6.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.374'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.99%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '148.43'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '174.19'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.12%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '43.24'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.33%  '
$ws.Range('E48').Value = '  +3.84%  '
$ws.Range('E49').Value = '  -9.42%  '
$ws.Range('E50').Value = '  -7.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.89'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.19%  '
